$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a cell's value as plain text (avoids Excel's automatic
# date-string-to-date-serial conversion), while leaving the cell's style
# as the default "Normal" style (no explicit number format applied).
function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Update date strings from DD/MM/YYYY to DD-MM-YYYY format (rows 3-21, column A)
Set-TextValue "A3" "28-07-2022"
Set-TextValue "A4" "01-08-2022"
Set-TextValue "A5" "04-08-2022"
Set-TextValue "A6" "08-08-2022"
Set-TextValue "A7" "11-08-2022"
Set-TextValue "A8" "15-08-2022"
Set-TextValue "A9" "18-08-2022"
Set-TextValue "A10" "22-08-2022"
Set-TextValue "A11" "25-08-2022"
Set-TextValue "A12" "29-08-2022"
Set-TextValue "A13" "01-09-2022"
Set-TextValue "A14" "05-09-2022"
Set-TextValue "A15" "08-09-2022"
Set-TextValue "A16" "12-09-2022"
Set-TextValue "A17" "15-09-2022"
Set-TextValue "A18" "19-09-2022"
Set-TextValue "A19" "22-09-2022"
Set-TextValue "A20" "26-09-2022"
Set-TextValue "A21" "29-09-2022"

# Update attendance counts for row 3 (28-07-2022)
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Update attendance counts for row 4 (01-08-2022)
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# Update attendance counts for row 5 (04-08-2022)
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

# Update attendance counts for row 11 (25-08-2022)
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("H11").Value = 0
